# Swap the roles of the first two sheets:
#  - the (empty) second sheet becomes the first sheet, renamed "titles",
#    and gets populated with a title/level table of contents.
#  - the first sheet (website/vulnerability data) becomes the second
#    sheet, renamed "vuls".
# Sheet3 is left untouched.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)
$sheet2 = $wb.Worksheets.Item(2)

# Moving (rather than just renaming) keeps each worksheet's identity -
# and therefore its sheetId / relationship id - attached to its data as
# it changes tab position.
$sheet2.Move($sheet1)

$wb.Worksheets.Item(1).Name = "titles"
$wb.Worksheets.Item(2).Name = "vuls"

# Renaming a sheet clears the sheet-qualifier on defined names that
# pointed at it, so restore the "vuls!" prefix on the old ToC refs.
foreach ($definedName in $wb.Names) {
    $definedName.RefersTo = "=vuls!#REF!"
}

# Populate the new "titles" sheet with a text/level outline.
$titles = $wb.Worksheets.Item(1)
$titles.Range("A2").Value = "title1"
$titles.Range("B1").Value = "level"
$titles.Range("A1").Value = "text"
$titles.Range("B2").Value = 1
$titles.Range("A3").Value = "title1.1"
$titles.Range("B3").Value = 2
$titles.Range("A4").Value = "title1.1.1"
$titles.Range("B4").Value = 3
$titles.Range("A5").Value = "title2"
$titles.Range("B5").Value = 1
$titles.Range("A6").Value = "title2.1"
$titles.Range("B6").Value = 2
$titles.Range("A7").Value = "title3"
$titles.Range("B7").Value = 1
$titles.Range("A8").Value = "title3.1"
$titles.Range("B8").Value = 2
$titles.Range("A9").Value = "title3.1.1"
$titles.Range("B9").Value = 3

[void]$titles.Range("A21").Select()
